$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 2000
$ws.Range("L2").Interior.Color = 65535
